# OY-3070 validate presence of hakemus OID when henkilo OID present
#
# Extends the "erillishaku_oidilla" test sheet with six new columns
# (Kutsumanimi, Syntymäpaikka, Passin numero, Kansallinen ID-tunnus,
# Kaupunki ja maa, Hakemus-oid) and renames the sample hakija-oid value
# so it reads as an oid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the used range (AE:AJ, rows 1-10) by replicating the
#    existing last-column formatting (style 2 on rows 1-2, style 7 on
#    rows 3-10) from column AD so the new cells look like the rest of
#    the table.
$ws.Range("AD1:AD10").Copy($ws.Range("AE1:AJ10"))

# 2) New header row (row 1) labels for the six new columns.
$ws.Range("AE1").Value = "Kutsumanimi"
$ws.Range("AF1").Value = "Syntymäpaikka"
$ws.Range("AG1").Value = "Passin numero"
$ws.Range("AH1").Value = "Kansallinen ID-tunnus"
$ws.Range("AI1").Value = "Kaupunki ja maa"
$ws.Range("AJ1").Value = "Hakemus-oid"

# 3) Sample data row (row 2) values for the new columns. AF2-AI2 stay
#    blank (no sample data supplied for those fields).
$ws.Range("AE2").Value = "Tuomas"
$ws.Range("AJ2").Value = "Hakemus1"

# 4) The applicant identifier sample value is renamed from the
#    placeholder "hakija1" to "Hakijaoid1" so it reads as an oid.
$ws.Range("G2").Value = "Hakijaoid1"
